$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 248, shifting existing rows 248:270 down to 249:271
$ws.Rows.Item(248).Insert()

# Populate the new row 248 with the new record's data.
# Unchanged-format columns (A,B,C,E,F,G,H,I,N,Q,R) are copied from the
# pattern shared by the whole block; only D,J,K,L,M,O,P are new values.
$ws.Cells.Item(248, 1).Value = 5
$ws.Cells.Item(248, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(248, 3).Value = "Maule"
$ws.Cells.Item(248, 4).Value = 44578
$ws.Cells.Item(248, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(248, 5).Value = 7
$ws.Cells.Item(248, 6).Value = 100112032
$ws.Cells.Item(248, 7).Value = "Zapallo italiano"
$ws.Cells.Item(248, 8).Value = "Sin especificar"
$ws.Cells.Item(248, 9).Value = "Primera"
$ws.Cells.Item(248, 10).Value = 400
$ws.Cells.Item(248, 11).Value = 7000
$ws.Cells.Item(248, 12).Value = 7000
$ws.Cells.Item(248, 13).Value = 7000
$ws.Cells.Item(248, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(248, 15).Value = "Región del Maule"
$ws.Cells.Item(248, 16).Value = 117
$ws.Cells.Item(248, 17).Value = 60
$ws.Cells.Item(248, 18).Value = "Hortaliza"
